$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.048025250434875
$ws.Range("B1").Value = 3.099463701248169
$ws.Range("C1").Value = 6.689461708068848
$ws.Range("D1").Value = 1.872572541236877
$ws.Range("E1").Value = 1.307311415672302
